# Append 12 new apartment-complex rows to the tracking list on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(101273, "자연앤힐스테이트"),
    @(111038, "광교중흥에스클래스(주상복합)"),
    @(102119, "광교e편한세상2차"),
    @(109412, "힐스테이트영통"),
    @(103518, "래미안영통마크원2단지"),
    @(105153, "한양수자인에듀파크"),
    @(109929, "영통라온프라이빗"),
    @(1804,   "영통에듀파크"),
    @(2490,   "건영1차"),
    @(22911,  "매탄위브하늘채"),
    @(135404, "영흥숲푸르지오파크비엔"),
    @(13698,  "현대힐스테이트")
)

$startRow = 205
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Move the selection / view down to the newly-added rows, matching the
# author's on-screen state after the paste (selection A205:A216).
$lastRow = $startRow + $newRows.Count - 1
$selRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($lastRow, 1))
$excel.Goto($selRange)
